$d = $word.ActiveDocument

# Mapping of original (Bosnian/Croatian/Serbian) labels to their English
# translations. Each paragraph currently holds a single run whose text is
# "<OldLabel>: ${placeholder}". The target shape keeps the same visible
# text but splits it into two runs: one that holds just the (now English)
# label, and one that holds the ": ${placeholder}" remainder.
$labels = @{
    'Primalac' = 'Receiver'
    'Platioc'  = 'Payer'
    'Usluge'   = 'Services'
    'Cijena'   = 'Cost'
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $pRange = $para.Range
    $text = $pRange.Text

    foreach ($old in $labels.Keys) {
        $prefix = $old + ':'
        if ($text.StartsWith($prefix)) {
            $new = $labels[$old]

            $paraStart = $pRange.Start
            $paraEnd = $pRange.End

            # Trim the trailing paragraph mark from the end offset so the
            # second run does not swallow it.
            $contentEnd = $paraEnd - 1

            # Range covering just the old label text (e.g. "Primalac").
            $labelRange = $d.Range($paraStart, $paraStart + $old.Length)
            $labelRange.Text = $new

            # The label text length may have changed; the remainder now
            # starts right after the (new) label and runs to the end of
            # the paragraph's visible content (": `${placeholder}").
            $newLabelEnd = $paraStart + $new.Length
            $remainderRange = $d.Range($newLabelEnd, $contentEnd - $old.Length + $new.Length)

            # Forcing a (no-op) formatting toggle on the remainder is what
            # makes the engine actually split it into its own run instead
            # of silently re-merging it with the identically formatted
            # label run.
            $remainderRange.Font.Bold = $true
            $remainderRange.Font.Bold = $false

            break
        }
    }
}
